# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Coin names/links/prices/%-volume are stored as plain text in the sheet, so
# every write is forced to Text first (NumberFormat "@") before assigning the
# value - this stops Excel from auto-coercing numeric-looking strings like
# "6.41" or "574.66" into real numbers. The cell's style is then restored to
# "Normal" so we don't leave a stray text-format style applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '66.608.78'
Set-TextValue 'E2' '  -0.33%  '
Set-TextValue 'D3' '3.070.92'
Set-TextValue 'E3' '  -1.35%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '574.66'
Set-TextValue 'E5' '  -0.62%  '
Set-TextValue 'D6' '168.97'
Set-TextValue 'E6' '  -1.72%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '3.067.87'
Set-TextValue 'E8' '  -1.34%  '
Set-TextValue 'E9' '  -2.13%  '
Set-TextValue 'D10' '6.41'
Set-TextValue 'E10' '  -0.60%  '
Set-TextValue 'E11' '  -1.83%  '
Set-TextValue 'E12' '  -3.44%  '
Set-TextValue 'E13' '  -2.74%  '
Set-TextValue 'D14' '35.54'
Set-TextValue 'E14' '  -4.11%  '
Set-TextValue 'E15' '  -1.78%  '
Set-TextValue 'D16' '3.582.31'
Set-TextValue 'E16' '  -1.32%  '
Set-TextValue 'D17' '66.556.77'
Set-TextValue 'E17' '  -0.36%  '
Set-TextValue 'E18' '  +3.88%  '
Set-TextValue 'D19' '6.94'
Set-TextValue 'E19' '  -3.34%  '
Set-TextValue 'D20' '3.071.46'
Set-TextValue 'E20' '  -1.34%  '
Set-TextValue 'D21' '486.93'
Set-TextValue 'E21' '  +2.22%  '
Set-TextValue 'E22' '  -2.50%  '
Set-TextValue 'E23' '  -3.90%  '
Set-TextValue 'E24' '  -1.85%  '
Set-TextValue 'E25' '  -4.73%  '
Set-TextValue 'E26' '  -3.54%  '
Set-TextValue 'D27' '10.13'
Set-TextValue 'E27' '  -0.54%  '
Set-TextValue 'E28' '  +0.07%  '
Set-TextValue 'D29' '7.78'
Set-TextValue 'E29' '  -1.13%  '
Set-TextValue 'E30' '  -4.95%  '
Set-TextValue 'E31' '  -2.60%  '
Set-TextValue 'E32' '  -3.83%  '
Set-TextValue 'E33' '  -3.60%  '
Set-TextValue 'E34' '  -3.61%  '
Set-TextValue 'E35' '  +0.03%  '
Set-TextValue 'B36' 'Mantle'
Set-TextValue 'C36' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D36' '0.947'
Set-TextValue 'E36' '  -2.74%  '
Set-TextValue 'B37' 'Filecoin'
Set-TextValue 'C37' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D37' '5.58'
Set-TextValue 'E37' '  -4.86%  '
Set-TextValue 'D38' '47.00'
Set-TextValue 'E38' '  -0.02%  '
Set-TextValue 'E39' '  -0.63%  '
Set-TextValue 'E40' '  -4.99%  '
Set-TextValue 'E41' '  -3.41%  '
Set-TextValue 'E42' '  -4.85%  '
Set-TextValue 'D43' '2.761.63'
Set-TextValue 'E43' '  -2.48%  '
Set-TextValue 'E44' '  -2.49%  '
Set-TextValue 'E45' '  -3.39%  '
Set-TextValue 'D46' '135.18'
Set-TextValue 'E46' '  -0.12%  '
Set-TextValue 'D47' '365.76'
Set-TextValue 'E47' '  -5.60%  '
Set-TextValue 'E48' '  +0.00%  '
Set-TextValue 'D49' '24.54'
Set-TextValue 'E49' '  -0.86%  '
Set-TextValue 'E50' '  -2.24%  '
Set-TextValue 'E51' '  -1.97%  '
